$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddVacancy")

# Update the existing row 2 "VacancyName" value from "accounts" to "test1"
$ws.Range("B2").Value = "test1"

# Add two more vacancy rows, mirroring row 2's other columns, to exercise
# "get no of columns" style tests across multiple rows.
$ws.Range("A3").Value = $ws.Range("A2").Text
$ws.Range("B3").Value = "test2"
$ws.Range("C3").Value = $ws.Range("C2").Text
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = $ws.Range("E2").Text

$ws.Range("A4").Value = $ws.Range("A2").Text
$ws.Range("B4").Value = "test3"
$ws.Range("C4").Value = $ws.Range("C2").Text
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = $ws.Range("E2").Text

$ws.Range("F4").Select()
